$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so values like "1.000"
# or "25.806.15" are not reinterpreted as numbers/dates, matching the
# original inline-string cell type.
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '25.806.15'
$ws.Range('E2').Value = '  -3.79%  '
$ws.Range('D3').Value = '1.815.54'
$ws.Range('E3').Value = '  -3.11%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = '277.30'
$ws.Range('E5').Value = '  -7.92%  '
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').Value = '  -0.11%  '
$ws.Range('D7').Value = '0.5098'
$ws.Range('E7').Value = '  -4.59%  '
$ws.Range('D8').Value = '0.3528'
$ws.Range('E8').Value = '  -6.16%  '
$ws.Range('D9').Value = '44.51'
$ws.Range('E9').Value = '  -2.04%  '
$ws.Range('D10').Value = '0.06669'
$ws.Range('E10').Value = '  -7.11%  '
$ws.Range('D11').Value = '20.03'
$ws.Range('E11').Value = '  -7.29%  '
$ws.Range('D12').Value = '0.8290'
$ws.Range('E12').Value = '  -6.60%  '
$ws.Range('E13').Value = '  -3.87%  '
$ws.Range('D14').Value = '1.815.83'
$ws.Range('E14').Value = '  -0.76%  '
$ws.Range('D15').Value = '5.077'
$ws.Range('E15').Value = '  -3.85%  '
$ws.Range('D16').Value = '87.49'
$ws.Range('E16').Value = '  -6.55%  '
$ws.Range('D17').Value = '0.9999'
$ws.Range('E17').Value = '  -0.07%  '
$ws.Range('D18').Value = '14.14'
$ws.Range('E18').Value = '  -4.30%  '
$ws.Range('D19').Value = '0.000008039'
$ws.Range('E19').Value = '  -5.98%  '
$ws.Range('E20').Value = '  -0.09%  '
$ws.Range('D21').Value = '25.850.41'
$ws.Range('E21').Value = '  -3.75%  '
$ws.Range('D22').Value = '4.731'
$ws.Range('E22').Value = '  -5.09%  '
$ws.Range('D23').Value = '10.02'
$ws.Range('E23').Value = '  -6.21%  '
$ws.Range('D24').Value = '6.086'
$ws.Range('E24').Value = '  -4.86%  '
$ws.Range('D25').Value = '141.37'
$ws.Range('E25').Value = '  -3.46%  '
$ws.Range('E26').Value = '  -3.38%  '
$ws.Range('D27').Value = '1.673'
$ws.Range('E27').Value = '  -3.88%  '
$ws.Range('D28').Value = '17.07'
$ws.Range('E28').Value = '  -5.15%  '
$ws.Range('D29').Value = '109.33'
$ws.Range('E29').Value = '  -4.01%  '
$ws.Range('D30').Value = '4.353'
$ws.Range('E30').Value = '  -7.82%  '
$ws.Range('D31').Value = '4.240'
$ws.Range('E31').Value = '  -7.93%  '
$ws.Range('D32').Value = '0.08801'
$ws.Range('E32').Value = '  -3.86%  '
$ws.Range('D33').Value = '0.04912'
$ws.Range('E33').Value = '  -1.21%  '
$ws.Range('D34').Value = '0.7272'
$ws.Range('E34').Value = '  -9.66%  '
$ws.Range('D35').Value = '1.138'
$ws.Range('E35').Value = '  -3.11%  '
$ws.Range('D36').Value = '2.876'
$ws.Range('E36').Value = '  -3.25%  '
$ws.Range('D37').Value = '1.000'
$ws.Range('E37').Value = '  -0.30%  '
$ws.Range('D38').Value = '3.132'
$ws.Range('E38').Value = '  -2.72%  '
$ws.Range('E39').Value = '  -8.05%  '
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').Value = '0.01853'
$ws.Range('E40').Value = '  -5.04%  '
$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D41').Value = '0.5191'
$ws.Range('E41').Value = '  -13.54%  '
$ws.Range('D42').Value = '0.9556'
$ws.Range('E42').Value = '  -10.78%  '
$ws.Range('D43').Value = '6.206'
$ws.Range('E43').Value = '  -5.10%  '
$ws.Range('D44').Value = '111.48'
$ws.Range('E44').Value = '  -3.00%  '
$ws.Range('D45').Value = '8.008'
$ws.Range('E45').Value = '  -9.77%  '
$ws.Range('D47').Value = '0.4570'
$ws.Range('E47').Value = '  -10.60%  '
$ws.Range('D48').Value = '0.1365'
$ws.Range('E48').Value = '  -8.37%  '
$ws.Range('D49').Value = '36.68'
$ws.Range('E49').Value = '  -2.29%  '
$ws.Range('D50').Value = '9.239'
$ws.Range('E50').Value = '  -6.58%  '
$ws.Range('D51').Value = '1.503'
$ws.Range('E51').Value = '  -7.80%  '

# Restore the default (unstyled) cell style now that the text values are
# committed, so no stray number-format style is left attached to the cells.
$ws.Range('D2:D51').Style = 'Normal'
